$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl
$ws.Cells.Item(18, 8).Value = 2750
$ws.Cells.Item(18, 9).Value = 500
$ws.Cells.Item(18, 11).Value = 500
$ws.Cells.Item(18, 13).Value = -216
# Row 40: Stuck in the Moment
$ws.Cells.Item(40, 8).Value = 4851.7407
$ws.Cells.Item(40, 9).Value = 1749.9333
$ws.Cells.Item(40, 10).Value = 8729
$ws.Cells.Item(40, 11).Value = 1749.9333
$ws.Cells.Item(40, 12).Value = 8729
$ws.Cells.Item(40, 13).Value = -1574.9333
$ws.Cells.Item(40, 14).Value = -9079
# Row 53: No Accounting for Waste
$ws.Cells.Item(53, 8).Value = 634.2
$ws.Cells.Item(53, 9).Value = 417.75
$ws.Cells.Item(53, 10).Value = 1500
$ws.Cells.Item(53, 11).Value = 417.75
$ws.Cells.Item(53, 12).Value = 1500
$ws.Cells.Item(53, 13).Value = 219.25
$ws.Cells.Item(53, 14).Value = -2774
# Row 96: Scroll Down
$ws.Cells.Item(96, 8).Value = 445.4
$ws.Cells.Item(96, 9).Value = 272.66666
$ws.Cells.Item(96, 10).Value = 2000
$ws.Cells.Item(96, 11).Value = 817.9999799999999
$ws.Cells.Item(96, 12).Value = 6000
$ws.Cells.Item(96, 13).Value = 555.0000200000001
$ws.Cells.Item(96, 14).Value = -8746
# Row 107: Another Man's Ink
$ws.Cells.Item(107, 8).Value = 1610
$ws.Cells.Item(107, 9).Value = 950.125
$ws.Cells.Item(107, 11).Value = 950.125
$ws.Cells.Item(107, 13).Value = 969.875
# Row 116: Growing Up
$ws.Cells.Item(116, 8).Value = 3999.5
$ws.Cells.Item(116, 9).Value = 3999.5
$ws.Cells.Item(116, 11).Value = 3999.5
$ws.Cells.Item(116, 13).Value = -557.5
# Row 135: For Tired Minds
$ws.Cells.Item(135, 8).Value = 1326.8667
$ws.Cells.Item(135, 10).Value = 3109.3333
$ws.Cells.Item(135, 12).Value = 27983.9997
$ws.Cells.Item(135, 14).Value = -33053.9997
# Row 137: Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 3723.6667
$ws.Cells.Item(137, 9).Value = 3388.0833
$ws.Cells.Item(137, 10).Value = 3947.389
$ws.Cells.Item(137, 11).Value = 10164.2499
$ws.Cells.Item(137, 12).Value = 11842.167
$ws.Cells.Item(137, 13).Value = -7614.249899999999
$ws.Cells.Item(137, 14).Value = -16942.167

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 122: Haste for High Durium
$ws.Cells.Item(122, 8).Value = 2875
$ws.Cells.Item(122, 9).Value = 2666.6667
$ws.Cells.Item(122, 11).Value = 8000.000100000001
$ws.Cells.Item(122, 13).Value = -5550.000100000001
# Row 132: Don't Bore Me, Ore Me
$ws.Cells.Item(132, 8).Value = 1998.5
$ws.Cells.Item(132, 9).Value = 1998.5
$ws.Cells.Item(132, 11).Value = 5995.5
$ws.Cells.Item(132, 13).Value = -3465.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 5: Axe Me Anything
$ws.Cells.Item(5, 8).Value = 1048.25
$ws.Cells.Item(5, 9).Value = 1048.25
$ws.Cells.Item(5, 11).Value = 1048.25
$ws.Cells.Item(5, 13).Value = -935.25
# Row 86: Through Thick and Thin
$ws.Cells.Item(86, 8).Value = 3279.8
$ws.Cells.Item(86, 9).Value = 3349.75
$ws.Cells.Item(86, 10).Value = 3000
$ws.Cells.Item(86, 11).Value = 3349.75
$ws.Cells.Item(86, 12).Value = 3000
$ws.Cells.Item(86, 13).Value = -2226.75
$ws.Cells.Item(86, 14).Value = -5246
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Cells.Item(89, 8).Value = 3279.8
$ws.Cells.Item(89, 9).Value = 3349.75
$ws.Cells.Item(89, 10).Value = 3000
$ws.Cells.Item(89, 11).Value = 16748.75
$ws.Cells.Item(89, 12).Value = 15000
$ws.Cells.Item(89, 13).Value = -11132.75
$ws.Cells.Item(89, 14).Value = -26232
# Row 94: High Steal
$ws.Cells.Item(94, 8).Value = 2197.125
$ws.Cells.Item(94, 9).Value = 2263
$ws.Cells.Item(94, 10).Value = 1999.5
$ws.Cells.Item(94, 11).Value = 2263
$ws.Cells.Item(94, 12).Value = 1999.5
$ws.Cells.Item(94, 13).Value = -1812
$ws.Cells.Item(94, 14).Value = -2901.5
# Row 105: Ingot to Wing It
$ws.Cells.Item(105, 8).Value = 2724.2222
$ws.Cells.Item(105, 9).Value = 2586.3333
$ws.Cells.Item(105, 11).Value = 2586.3333
$ws.Cells.Item(105, 13).Value = -839.3332999999998
# Row 134: Ruthenium Supremium
$ws.Cells.Item(134, 8).Value = 9899.200000000001
$ws.Cells.Item(134, 10).Value = 16499
$ws.Cells.Item(134, 12).Value = 49497
$ws.Cells.Item(134, 14).Value = -54567

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent
$ws.Cells.Item(7, 8).Value = 765
$ws.Cells.Item(7, 9).Value = 647.5
$ws.Cells.Item(7, 10).Value = 1000
$ws.Cells.Item(7, 11).Value = 647.5
$ws.Cells.Item(7, 12).Value = 1000
$ws.Cells.Item(7, 13).Value = -534.5
$ws.Cells.Item(7, 14).Value = -1226
# Row 99: O Pine
$ws.Cells.Item(99, 8).Value = 3999
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 3999
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 3999
$ws.Cells.Item(99, 13).ClearContents()
$ws.Cells.Item(99, 14).Value = -6995
# Row 122: Timber of Tenkonto
$ws.Cells.Item(122, 8).Value = 2999.6
$ws.Cells.Item(122, 9).Value = 2999.75
$ws.Cells.Item(122, 10).Value = 2999
$ws.Cells.Item(122, 11).Value = 8999.25
$ws.Cells.Item(122, 12).Value = 8997
$ws.Cells.Item(122, 13).Value = -6549.25
$ws.Cells.Item(122, 14).Value = -13897
# Row 126: A Better Conductor
$ws.Cells.Item(126, 8).Value = 3999
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 3999
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 11997
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 14).Value = -16937

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 34: Fever Pitch
$ws.Cells.Item(34, 8).Value = 2183.1428
$ws.Cells.Item(34, 10).Value = 2530.3333
$ws.Cells.Item(34, 12).Value = 7590.999899999999
$ws.Cells.Item(34, 14).Value = -7758.999899999999
# Row 63: The Next to Last Supper
$ws.Cells.Item(63, 8).Value = 829
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 14).ClearContents()
# Row 66: Nostalgia through the Stomach (L)
$ws.Cells.Item(66, 8).Value = 829
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 14).ClearContents()
# Row 137: Creative Chocolate
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).ClearContents()
$ws.Cells.Item(137, 14).ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43: Get the Green Stuff
$ws.Cells.Item(43, 8).Value = 9322.727999999999
$ws.Cells.Item(43, 10).Value = 11428.571
$ws.Cells.Item(43, 12).Value = 11428.571
$ws.Cells.Item(43, 14).Value = -11730.571
# Row 126: Gold Rush Order
$ws.Cells.Item(126, 8).Value = 5149
$ws.Cells.Item(126, 9).Value = 5635.6
$ws.Cells.Item(126, 11).Value = 16906.8
$ws.Cells.Item(126, 13).Value = -14436.8
# Row 132: On Board for Lar
$ws.Cells.Item(132, 8).Value = 6333.3335
$ws.Cells.Item(132, 9).Value = 5000
$ws.Cells.Item(132, 10).Value = 7000
$ws.Cells.Item(132, 11).Value = 15000
$ws.Cells.Item(132, 12).Value = 21000
$ws.Cells.Item(132, 13).Value = -12470
$ws.Cells.Item(132, 14).Value = -26060
# Row 136: Shiny and Good
$ws.Cells.Item(136, 8).Value = 80000
$ws.Cells.Item(136, 10).Value = 80000
$ws.Cells.Item(136, 12).Value = 240000
$ws.Cells.Item(136, 14).Value = -245100

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Cells.Item(7, 8).Value = 5603.6
$ws.Cells.Item(7, 9).Value = 5004.5
$ws.Cells.Item(7, 10).Value = 8000
$ws.Cells.Item(7, 11).Value = 5004.5
$ws.Cells.Item(7, 12).Value = 8000
$ws.Cells.Item(7, 13).Value = -4892.5
$ws.Cells.Item(7, 14).Value = -8224
# Row 18: Simply the Best
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).ClearContents()
# Row 20: Choke Hold
$ws.Cells.Item(20, 8).Value = 20000
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 13).ClearContents()
# Row 22: Skin off Their Backs
$ws.Cells.Item(22, 8).Value = 850
# Row 27: Fire and Hide
$ws.Cells.Item(27, 8).Value = 850
# Row 61: Spelling Me Softly
$ws.Cells.Item(61, 8).Value = 7198.4
$ws.Cells.Item(61, 9).Value = 6998
$ws.Cells.Item(61, 11).Value = 6998
$ws.Cells.Item(61, 13).Value = -6796
# Row 93: Hide to Go Seek
$ws.Cells.Item(93, 8).Value = 1728.1538
$ws.Cells.Item(93, 9).Value = 1830.5
$ws.Cells.Item(93, 10).Value = 500
$ws.Cells.Item(93, 11).Value = 1830.5
$ws.Cells.Item(93, 12).Value = 500
$ws.Cells.Item(93, 13).Value = -582.5
$ws.Cells.Item(93, 14).Value = -2996
# Row 100: Tiger in the Sack
$ws.Cells.Item(100, 8).Value = 3310.4
$ws.Cells.Item(100, 9).Value = 3178.2222
$ws.Cells.Item(100, 10).Value = 4500
$ws.Cells.Item(100, 11).Value = 3178.2222
$ws.Cells.Item(100, 12).Value = 4500
$ws.Cells.Item(100, 13).Value = -2637.2222
$ws.Cells.Item(100, 14).Value = -5582
# Row 113: Peace in Rest
$ws.Cells.Item(113, 8).Value = 7198.4
$ws.Cells.Item(113, 9).Value = 6998
$ws.Cells.Item(113, 11).Value = 6998
$ws.Cells.Item(113, 13).Value = -4828
# Row 122: Hell on Leather
$ws.Cells.Item(122, 8).Value = 8000
$ws.Cells.Item(122, 9).Value = 8000
$ws.Cells.Item(122, 11).Value = 24000
$ws.Cells.Item(122, 13).Value = -21550
# Row 126: Battered Books
$ws.Cells.Item(126, 8).Value = 5603.6
$ws.Cells.Item(126, 9).Value = 5004.5
$ws.Cells.Item(126, 10).Value = 8000
$ws.Cells.Item(126, 11).Value = 15013.5
$ws.Cells.Item(126, 12).Value = 24000
$ws.Cells.Item(126, 13).Value = -12543.5
$ws.Cells.Item(126, 14).Value = -28940
# Row 132: Tenets of Tanning
$ws.Cells.Item(132, 8).Value = 24049.9
$ws.Cells.Item(132, 9).Value = 24785.572
$ws.Cells.Item(132, 10).Value = 22333.334
$ws.Cells.Item(132, 11).Value = 74356.716
$ws.Cells.Item(132, 12).Value = 67000.00199999999
$ws.Cells.Item(132, 13).Value = -71826.716
$ws.Cells.Item(132, 14).Value = -72060.00199999999
# Row 136: Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 12497.5
$ws.Cells.Item(136, 10).Value = 6500
$ws.Cells.Item(136, 12).Value = 19500
$ws.Cells.Item(136, 14).Value = -24600

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 94: Proper Props
$ws.Cells.Item(94, 8).Value = 37500
$ws.Cells.Item(94, 10).Value = 37500
$ws.Cells.Item(94, 12).Value = 37500
$ws.Cells.Item(94, 14).Value = -39302
# Row 100: Of Great Import
$ws.Cells.Item(100, 8).Value = 4875.0586
$ws.Cells.Item(100, 9).Value = 6077.231
$ws.Cells.Item(100, 10).Value = 968
$ws.Cells.Item(100, 11).Value = 12154.462
$ws.Cells.Item(100, 12).Value = 1936
$ws.Cells.Item(100, 13).Value = -11613.462
$ws.Cells.Item(100, 14).Value = -3018
# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 10399.4
$ws.Cells.Item(132, 9).Value = 7499.5
$ws.Cells.Item(132, 10).Value = 12332.667
$ws.Cells.Item(132, 11).Value = 22498.5
$ws.Cells.Item(132, 12).Value = 36998.001
$ws.Cells.Item(132, 13).Value = -19968.5
$ws.Cells.Item(132, 14).Value = -42058.001
# Row 136: Weaving the Envelope
$ws.Cells.Item(136, 8).Value = 1765.6666
$ws.Cells.Item(136, 9).Value = 1765.6666
$ws.Cells.Item(136, 11).Value = 5296.9998
$ws.Cells.Item(136, 13).Value = -2746.9998
